$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff reflects a cyclic swap of the data rows (2<->4 and 3<->5) across
# columns D (Fecha) and K..T (Variedad .. Kg/unidad). Columns A,B,C,E..J stay
# the same since they hold values shared across the rows being swapped.

# Row 2 new values (was row 4)
$ws.Range("D2").Value = 44181
$ws.Range("K2").Value = "Modesto"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("Q2").Value = "$/caja 18 kilos"
$ws.Range("R2").Value = "Región de Coquimbo"
$ws.Range("S2").Value = 1139
$ws.Range("T2").Value = 18

# Row 3 new values (was row 5)
$ws.Range("D3").Value = 44174
$ws.Range("K3").Value = "Castle Brite"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 75
$ws.Range("N3").Value = 9000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 9467
$ws.Range("Q3").Value = "$/caja 10 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 947
$ws.Range("T3").Value = 10

# Row 4 new values (was row 3)
$ws.Range("D4").Value = 44189
$ws.Range("K4").Value = "Dina"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 16562
$ws.Range("Q4").Value = "$/caja 18 kilos"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 920
$ws.Range("T4").Value = 18

# Row 5 new values (was row 2)
$ws.Range("D5").Value = 44165
$ws.Range("K5").Value = "Castle Brite"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 17000
$ws.Range("P5").Value = 16500
$ws.Range("Q5").Value = "$/caja 15 kilos granel"
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 1100
$ws.Range("T5").Value = 15
